# Applies the cryptos.xlsx price/volume refresh described in the commit
# message: "Updated cryptos list on Mon Apr 17 23:30:34 UTC 2023 with
# GitHub Actions" -- refreshed Price/Volume(1h) figures for every coin row,
# plus a rank swap between Polkadot/Chainlink (rows 14-15) and between
# TrustWalletToken/Algorand (rows 42-43).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.716.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.58%  "

# Row 3
$ws.Range("D3").Value = "'2.095.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.81%  "

# Row 4
$ws.Range("D4").Value = "'1.011"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.25%  "

# Row 5
$ws.Range("D5").Value = "'343.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.28%  "

# Row 6
$ws.Range("E6").Value = "  +0.14%  "

# Row 7
$ws.Range("D7").Value = "'0.5163"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.73%  "

# Row 8
$ws.Range("D8").Value = "'0.4381"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.95%  "

# Row 9
$ws.Range("D9").Value = "'52.99"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.90%  "

# Row 10
$ws.Range("D10").Value = "'0.09225"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.16%  "

# Row 11
$ws.Range("E11").Value = "  -1.99%  "

# Row 12
$ws.Range("D12").Value = "'24.88"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.47%  "

# Row 13
$ws.Range("D13").Value = "'2.102.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.67%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'8.244"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.84%  "

# Row 15
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'6.772"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.65%  "

# Row 16
$ws.Range("D16").Value = "'99.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.64%  "

# Row 17
$ws.Range("E17").Value = "  -1.35%  "

# Row 18
$ws.Range("D18").Value = "'1.010"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.11%  "

# Row 19
$ws.Range("D19").Value = "'20.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.96%  "

# Row 20
$ws.Range("D20").Value = "'0.06657"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.94%  "

# Row 21
$ws.Range("E21").Value = "  +0.16%  "

# Row 22
$ws.Range("D22").Value = "'6.206"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.65%  "

# Row 23
$ws.Range("D23").Value = "'29.753.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.80%  "

# Row 24
$ws.Range("E24").Value = "  -2.42%  "

# Row 25
$ws.Range("E25").Value = "  -2.60%  "

# Row 26
$ws.Range("D26").Value = "'2.344.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.87%  "

# Row 27
$ws.Range("D27").Value = "'21.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.29%  "

# Row 28
$ws.Range("E28").Value = "  -2.94%  "

# Row 29
$ws.Range("D29").Value = "'161.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.33%  "

# Row 30
$ws.Range("D30").Value = "'133.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.78%  "

# Row 31
$ws.Range("D31").Value = "'1.139"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.80%  "

# Row 32
$ws.Range("E32").Value = "  -2.93%  "

# Row 33
$ws.Range("D33").Value = "'1.636"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.55%  "

# Row 34
$ws.Range("D34").Value = "'6.174"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.38%  "

# Row 35
$ws.Range("D35").Value = "'3.950"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.11%  "

# Row 36
$ws.Range("D36").Value = "'6.303"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.67%  "

# Row 37
$ws.Range("D37").Value = "'10.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.83%  "

# Row 38
$ws.Range("D38").Value = "'0.02576"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.44%  "

# Row 39
$ws.Range("D39").Value = "'0.7098"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.75%  "

# Row 40
$ws.Range("D40").Value = "'0.06730"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.61%  "

# Row 41
$ws.Range("D41").Value = "'12.47"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.91%  "

# Row 42
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.2229"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.31%  "

# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'1.320"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.51%  "

# Row 44
$ws.Range("D44").Value = "'0.6989"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.44%  "

# Row 45
$ws.Range("D45").Value = "'14.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.57%  "

# Row 46
$ws.Range("E46").Value = "  +0.18%  "

# Row 47
$ws.Range("D47").Value = "'2.319"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.11%  "

# Row 48
$ws.Range("D48").Value = "'3.619"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.45%  "

# Row 50
$ws.Range("E50").Value = "  -2.30%  "

# Row 51
$ws.Range("D51").Value = "'82.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.92%  "
